{"js": "// Move the \"Mogu\u0107nost otkazivanja Premium paketa\" bullet from the\n// \"Registrovani korisnici\" list (numId=3) to the end of the\n// \"Premium korisnici\" list (numId=4), right after\n// \"Dobijanje notifikacije za po\u010detak/kraj utakmice\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the paragraph to move by its exact text.\nconst sourceText = \"Mogu\u0107nost otkazivanja Premium paketa\";\nlet sourceParagraph = null;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === sourceText) {\n    sourceParagraph = items[i];\n    break;\n  }\n}\nif (!sourceParagraph) {\n  throw new Error(\"Could not find paragraph: \" + sourceText);\n}\n\n// Locate the last item of the \"Premium korisnici\" list, i.e. the\n// paragraph right before the \"Administracija:\" heading.\nconst targetText = \"Dobijanje notifikacije za po\u010detak/kraj utakmice\";\nlet targetParagraph = null;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === targetText) {\n    targetParagraph = items[i];\n  }\n}\nif (!targetParagraph) {\n  throw new Error(\"Could not find paragraph: \" + targetText);\n}\n\n// Insert a new paragraph after the target (inherits its list/pPr\n// formatting automatically), carrying the moved sentence, then drop\n// the original paragraph from the \"Registrovani korisnici\" list.\ntargetParagraph.insertParagraph(sourceText, Word.InsertLocation.after);\nsourceParagraph.delete();\n\nawait context.sync();\n", "ps1": "# Move the \"Mogu\u0107nost otkazivanja Premium paketa\" bullet from the\n# \"Registrovani korisnici\" list (numId=3) to the end of the\n# \"Premium korisnici\" list (numId=4), right after\n# \"Dobijanje notifikacije za po\u010detak/kraj utakmice\".\n\n$d = $word.ActiveDocument\n\n$sourceText = \"Mogu\u0107nost otkazivanja Premium paketa\"\n$targetText = \"Dobijanje notifikacije za po\u010detak/kraj utakmice\"\n\n# Find the source paragraph (the one to move out of \"Registrovani korisnici\").\n$sourceParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.TrimEnd(\"`r\", \"`a\", \"`n\") -eq $sourceText) {\n        $sourceParagraph = $p\n        break\n    }\n}\nif ($null -eq $sourceParagraph) {\n    throw \"Could not find paragraph: $sourceText\"\n}\n\n# Find the last matching paragraph (the last item of \"Premium korisnici\",\n# right before the \"Administracija:\" heading).\n$targetParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.TrimEnd(\"`r\", \"`a\", \"`n\") -eq $targetText) {\n        $targetParagraph = $p\n    }\n}\nif ($null -eq $targetParagraph) {\n    throw \"Could not find paragraph: $targetText\"\n}\n\n# Insert a new paragraph right after the target (it inherits the target's\n# list/paragraph formatting, i.e. numId=4), set its text to the moved\n# sentence, then remove the original paragraph from the numId=3 list.\n$targetParagraph.Range.InsertParagraphAfter()\n$newParagraph = $targetParagraph.Next()\n$newParagraph.Range.Text = $sourceText\n\n$sourceParagraph.Range.Delete()\n"}
